$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spieltabelle")

# --- Row 16 (A16 = 14) ---
$ws.Range("B16").Value = "Cashgame"
$ws.Range("C16").Value = "sc.ch"
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0.8
$ws.Range("F16").Value = 0.01
$ws.Range("H16").Value = 45984
$ws.Range("I16").Value = 0.70833333333333337
$ws.Range("J16").Value = 45984
$ws.Range("K16").Value = 0.71527777777777779
$ws.Range("M16").Value = 1.2
$ws.Range("N16").Value = 0
$ws.Range("Q16").Value = 1.6
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 0.66
$ws.Range("T16").Value = "Nichts"
$ws.Range("U16").Value = "Ging All in Gegen tight aggresiv player mit TT. Er hatte JJ"
$ws.Range("V16").Value = "Bluffen auf Trockenen Felder mit einer hohen Karte funtkioniert erstaunlich gut"

# --- Row 17 (A17 = 15) ---
$ws.Range("B17").Value = "Cashgame"
$ws.Range("C17").Value = "sc.ch"
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0.8
$ws.Range("F17").Value = 0.01
$ws.Range("H17").Value = 45984
$ws.Range("I17").Value = 0.71527777777777779
$ws.Range("J17").Value = 45984
$ws.Range("K17").Value = 0.76736111111111116
$ws.Range("M17").Value = 1.2
$ws.Range("N17").Value = 0
$ws.Range("Q17").Value = 1.2
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 0.38
$ws.Range("T17").Value = "Nichts"
$ws.Range("V17").Value = "Bluff der Gegner richtig erkannt. Eigene haben  trotzdem funktioniert"
$ws.Range("U17").Value = "Ich wurde als Fisch erkannt und konnte dies teilweise nutzen, durchschnittliche Karten"

# --- Row 18 (A18 = 16) ---
$ws.Range("B18").Value = "Cashgame"
$ws.Range("C18").Value = "sc.ch"
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0.8
$ws.Range("F18").Value = 0.01
$ws.Range("H18").Value = 45984
$ws.Range("I18").Value = 0.76736111111111116
$ws.Range("J18").Value = 45984
$ws.Range("K18").Value = 0.77430555555555558
$ws.Range("M18").Value = 1.2
$ws.Range("N18").Value = 0
$ws.Range("Q18").Value = 1.2
$ws.Range("R18").Value = 3
$ws.Range("S18").Value = 0.62
$ws.Range("T18").Value = "Nichts"
$ws.Range("U18").Value = "Hatte Pech. Hatte überbettet, Gegner bekam trotzdem den Flush"
$ws.Range("V18").Value = "Verzweiflung führte Allin bei Mddlepair. Ging verloren"

# --- Row 19 (A19 = 17) ---
$ws.Range("B19").Value = "Cashgame"
$ws.Range("C19").Value = "sc.ch"
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0.8
$ws.Range("F19").Value = 0.01
$ws.Range("H19").Value = 45984
$ws.Range("I19").Value = 0.77430555555555558
$ws.Range("J19").Value = 45984
$ws.Range("K19").Value = 0.78402777777777777
$ws.Range("M19").Value = 1.2
$ws.Range("N19").Value = 0
$ws.Range("Q19").Value = 1.97
$ws.Range("R19").Value = 3
$ws.Range("S19").Value = 0.46
$ws.Range("T19").Value = "Nichts"
$ws.Range("U19").Value = "Hatte Pech. Hatte überbettet, ging mit 2 paaren all in, Gegner hatte Set getroffen"
$ws.Range("V19").Value = "Gegner waren gut lesbar"

# --- Row 20 (A20 = 18) : only a subset of columns filled in ---
$ws.Range("B20").Value = "Cashgame"
$ws.Range("C20").Value = "sc.ch"
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 0.8
$ws.Range("F20").Value = 0.01
$ws.Range("M20").Value = 1.2
$ws.Range("R20").Value = 3

# --- Selection / view state ---
$ws.Range("J21").Select()
